$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing row 3 IMEI value first (preserves shared-string order)
$ws.Range("B3").Value = "1594826001"

# Add a new row 4, mirroring rows 2/3's layout/format, with a new IMEI value.
# Copy row 3 (values+formats) down to row 4 so the text-typed columns (B/C/D)
# keep their text number format/style instead of being re-typed as numbers.
$ws.Range("A3:D3").Copy($ws.Range("A4:D4"))

$ws.Range("A4").Formula = "=FALSE()"
$ws.Range("B4").Value = "1594826003"
$ws.Range("C4").Value = "3.7.208.99:5175"
$ws.Range("D4").Value = "5"

# Move the active selection to F10
$ws.Range("F10").Select()
